$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-3 (values change only)
$ws.Range("C2").Value = 3.08
$ws.Range("C3").Value = 2.17

# Row 4: Cecilia Gigena removed -> becomes Celia del Carmen Dibernardi (was row5)
$ws.Range("A4").Value = "Celia del Carmen Dibernardi"
$ws.Range("B4").Value = 38
$ws.Range("C4").Value = 4.33

# Row 5: becomes Cristina Lorente (was row6)
$ws.Range("A5").Value = "Cristina Lorente "
$ws.Range("B5").Value = 47
$ws.Range("C5").Value = 5.36

# Row 6: becomes Florencia Marchese (was row7)
$ws.Range("A6").Value = "Florencia Marchese"
$ws.Range("B6").Value = 26
$ws.Range("C6").Value = 2.96

# Row 7: becomes Leticia Ainsa (was row8)
$ws.Range("A7").Value = "Leticia Ainsa"
$ws.Range("B7").Value = 22
$ws.Range("C7").Value = 2.51

# Row 8: becomes Liliana Bennice (was row9)
$ws.Range("A8").Value = "Liliana Bennice"
$ws.Range("B8").Value = 64
$ws.Range("C8").Value = 7.3

# Row 9: becomes MAGALI RIVAS (was row10)
$ws.Range("A9").Value = "MAGALI RIVAS"
$ws.Range("B9").Value = 36
$ws.Range("C9").Value = 4.1

# Row 10: becomes Magali Rivas (was row11, but different name+values)
$ws.Range("A10").Value = "Magali Rivas"
$ws.Range("B10").Value = 28
$ws.Range("C10").Value = 3.19

# Row 11: Marcela Ansonnaud stays but values change
$ws.Range("A11").Value = "Marcela Ansonnaud"
$ws.Range("B11").Value = 104
$ws.Range("C11").Value = 11.86

# Row 12: Marcia santa cruz -> Maria Teresa Onega
$ws.Range("A12").Value = "Maria Teresa Onega"
$ws.Range("B12").Value = 10
$ws.Range("C12").Value = 1.14

# Row 13: Mariana Sabbag stays, value changes
$ws.Range("A13").Value = "Mariana Sabbag"
$ws.Range("B13").Value = 57
$ws.Range("C13").Value = 6.5

# Row 14: Maximiliano Troncoso stays, values change
$ws.Range("A14").Value = "Maximiliano Troncoso"
$ws.Range("B14").Value = 75
$ws.Range("C14").Value = 8.55

# Row 15: Susana Lemmo -> Silvia Pisellini Marchegiani
$ws.Range("A15").Value = "Silvia Pisellini Marchegiani"
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = 0.23

# Row 16: Susana fernandez -> Susana Lemmo
$ws.Range("A16").Value = "Susana Lemmo"
$ws.Range("B16").Value = 322
$ws.Range("C16").Value = 36.72
